# Updated cryptos list on Mon May 22 13:23:08 UTC 2023 with GitHub Actions
# Refresh Coin/Link/Price/Volume(1h) columns (B:E) to match the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.975.26'
$ws.Range("E2").Value = '  -0.49%  '

# Row 3
$ws.Range("D3").Value = '1.827.40'
$ws.Range("E3").Value = '  +0.07%  '

# Row 4
$ws.Range("D4").Value = '''1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '''311.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.44%  '

# Row 6
$ws.Range("D6").Value = '''1.008'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '

# Row 7
$ws.Range("D7").Value = '''0.4657'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.62%  '

# Row 8
$ws.Range("D8").Value = '''0.3704'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.45%  '

# Row 9
$ws.Range("D9").Value = '''0.07366'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.31%  '

# Row 10
$ws.Range("D10").Value = '''0.8731'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.83%  '

# Row 11
$ws.Range("D11").Value = '''19.93'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.57%  '

# Row 12
$ws.Range("D12").Value = '''0.07826'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.56%  '

# Row 13
$ws.Range("D13").Value = '1.851.27'
$ws.Range("E13").Value = '  -2.04%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''5.342'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.52%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '''6.566'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.59%  '

# Row 16
$ws.Range("D16").Value = '''91.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.48%  '

# Row 17
$ws.Range("D17").Value = '''1.010'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.34%  '

# Row 18
$ws.Range("D18").Value = '''0.000008833'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.32%  '

# Row 19
$ws.Range("D19").Value = '''1.007'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '

# Row 20
$ws.Range("D20").Value = '''14.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.26%  '

# Row 21
$ws.Range("D21").Value = '26.679.91'
$ws.Range("E21").Value = '  -2.87%  '

# Row 22
$ws.Range("D22").Value = '''5.131'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.00%  '

# Row 23
$ws.Range("D23").Value = '''10.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.40%  '

# Row 24
$ws.Range("D24").Value = '2.011.71'
$ws.Range("E24").Value = '  -3.30%  '

# Row 25
$ws.Range("D25").Value = '''152.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.64%  '

# Row 26
$ws.Range("D26").Value = '''1.829'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.85%  '

# Row 27
$ws.Range("E27").Value = '  -1.13%  '

# Row 28
$ws.Range("D28").Value = '''2.078'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.87%  '

# Row 29
$ws.Range("D29").Value = '''5.106'
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = '''115.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.74%  '

# Row 31
$ws.Range("D31").Value = '''0.08878'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.53%  '

# Row 32
$ws.Range("D32").Value = '''2.963'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.68%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''4.435'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.68%  '

# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7249'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.47%  '

# Row 35
$ws.Range("D35").Value = '''1.133'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.54%  '

# Row 36
$ws.Range("D36").Value = '''2.480'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.18%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.01952'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.77%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''1.071'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.86%  '

# Row 39
$ws.Range("D39").Value = '''0.05214'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.24%  '

# Row 40
$ws.Range("D40").Value = '''2.920'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.43%  '

# Row 41
$ws.Range("D41").Value = '''7.150'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.52%  '

# Row 42
$ws.Range("D42").Value = '''0.5184'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.17%  '

# Row 43
$ws.Range("D43").Value = '''0.8670'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -13.94%  '

# Row 45
$ws.Range("D45").Value = '''8.200'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.19%  '

# Row 46
$ws.Range("D46").Value = '''0.4812'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.44%  '

# Row 47
$ws.Range("D47").Value = '''1.009'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.11%  '

# Row 48
$ws.Range("D48").Value = '''10.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.77%  '

# Row 49
$ws.Range("D49").Value = '''102.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.35%  '

# Row 50
$ws.Range("D50").Value = '''1.624'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '

# Row 51
$ws.Range("D51").Value = '''0.06209'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.93%  '
